$wb = $excel.ActiveWorkbook

# Add the new "strings" worksheet right after the existing "integer" sheet.
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "strings"

# Fill in the data in the same order it was originally typed (this controls
# the order new entries land in the shared-strings table).
$ws2.Range("A1").Value = "pass"
$ws2.Range("B1").Value = "W/C"

$ws2.Range("A2").Value = "pass"
$ws2.Range("B2").Value = "w/c"

$ws2.Range("A5").Value = "pass"
$ws2.Range("B5").Value = " W/C"

$ws2.Range("A6").Value = "pass"
$ws2.Range("B6").Value = " W/C "

$ws2.Range("A7").Value = "pass"
$ws2.Range("B7").Value = "W / C "

$ws2.Range("A3").Value = "fail"
$ws2.Range("B3").Value = "WC"

$ws2.Range("A8").Value = "pass"
$ws2.Range("B8").Value = "W/C #0000000000"

$ws2.Range("A4").Value = "fail"
$ws2.Range("B4").Value = " "

$ws2.Range("C4").Value = "<<B4 is a space"

# Column C was widened (e.g. to fit the "B4 is a space" note).
$ws2.Columns.Item(3).ColumnWidth = 12.83

# Leave the selection where the user last clicked, and make "strings" the
# active (selected) tab.
$ws2.Range("C10").Select()
$ws2.Activate()
